{"js": "// The document has a single \"_GoBack\" bookmark that Word keeps pinned to the\n// location of the most recent edit. The author deleted the single-letter \"I\"\n// run inside the `\"RI,R,0\\n\"` stream-protocol example text (in the\n// \"Stream protocol:\" paragraph whose reply string is `\"RI,R,0\\n\"`), which\n// caused Word to relocate the \"_GoBack\" bookmark from its previous spot\n// (end of the \"For I2C services...\" bullet) to the position where that\n// character was removed.\n\nconst doc = context.document;\nconst body = doc.body;\n\n// Remove the pre-existing \"_GoBack\" bookmark (only one may exist in a Word\n// document at a time); harmless no-op if it is not present.\ndoc.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// Find the unique paragraph text containing the \"RI,R,0\" example string.\nconst matches = body.search(\"RI,R,0\", { matchCase: true, matchWholeWord: false });\nmatches.load(\"items\");\nawait context.sync();\n\nif (matches.items.length === 0) {\n  throw new Error(\"Could not find the 'RI,R,0' example text to edit.\");\n}\n\n// Narrow down to just the single \"I\" character within that match.\nconst container = matches.items[0];\nconst iMatches = container.search(\"I\", { matchCase: true, matchWholeWord: false });\niMatches.load(\"items\");\nawait context.sync();\n\nif (iMatches.items.length === 0) {\n  throw new Error(\"Could not find the 'I' run to remove.\");\n}\n\nconst iRange = iMatches.items[0];\n\n// Insert the relocated \"_GoBack\" bookmark exactly where the \"I\" character\n// used to be, then delete that character (the order matters: inserting the\n// bookmark before deleting keeps it anchored to the collapsed position left\n// behind by the removed run, matching Word's own behavior).\niRange.insertBookmark(\"_GoBack\");\niRange.delete();\n\nawait context.sync();\n", "ps1": "# The document carries a single \"_GoBack\" bookmark that Word pins to the\n# location of the most recent edit. The author deleted the single-letter \"I\"\n# run inside the \"RI,R,0\\n\" stream-protocol example text (in the\n# \"Stream protocol:\" paragraph whose reply string is \"RI,R,0\\n\"), which moved\n# the \"_GoBack\" bookmark from its previous spot (end of the \"For I2C\n# services...\" bullet) to the position where that character was removed.\n\n$d = $word.ActiveDocument\n\n# Word keeps only one \"_GoBack\" bookmark in the whole document; remove the\n# old one so it does not linger once we add the new one below (no-op/safe if\n# it is already gone).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n  $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# Locate the unique \"RI,R,0\" example text using the document's plain-text\n# index (kept in sync with Range/Find character offsets).\n$fullText = $d.Content.Text\n$matchIndex = $fullText.IndexOf(\"RI,R,0\")\nif ($matchIndex -lt 0) {\n  throw \"Could not find the 'RI,R,0' example text to edit.\"\n}\n\n# The \"I\" is the character immediately after the \"R\" at the start of the\n# match; build a one-character Range addressing just that letter.\n$iRange = $d.Range($matchIndex + 1, $matchIndex + 2)\nif ($iRange.Text -ne \"I\") {\n  throw \"Unexpected text at the computed offset: '$($iRange.Text)'\"\n}\n\n# Delete the \"I\" character, then re-create the \"_GoBack\" bookmark collapsed\n# at the spot it used to occupy - matching Word's own \"last edit\" behavior.\n$iRange.Delete()\n$d.Bookmarks.Add(\"_GoBack\", $iRange)\n"}
